$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111609172
$ws.Range("B2").Value = 77268
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 228912
$ws.Range("F2").Value = 'Mörk kolflarnlav'
$ws.Range("G2").Value = 'Carbonicola myrmecina'
$ws.Range("H2").Value = '(Ach.) Bendiksby & Timdal'
$ws.Range("I2").Value = ""
$ws.Range("Q2").Value = 514955.9350709137
$ws.Range("R2").Value = 6925302.779521272
$ws.Range("AO2").Value = 'brandstubbe'

# Row 3
$ws.Range("A3").Value = 111609173
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = 'VU'
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = 'Knärot'
$ws.Range("G3").Value = 'Goodyera repens'
$ws.Range("H3").Value = '(L.) R. Br.'
$ws.Range("I3").Value = '''7'
$ws.Range("Q3").Value = 514934.1293421969
$ws.Range("R3").Value = 6925308.234934391
$ws.Range("AO3").Value = ""

# Row 4
$ws.Range("A4").Value = 111609175
$ws.Range("B4").Value = 77268
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 228912
$ws.Range("F4").Value = 'Mörk kolflarnlav'
$ws.Range("G4").Value = 'Carbonicola myrmecina'
$ws.Range("H4").Value = '(Ach.) Bendiksby & Timdal'
$ws.Range("I4").Value = ""
$ws.Range("Q4").Value = 514769.8196280882
$ws.Range("R4").Value = 6925156.6384242
$ws.Range("AO4").Value = 'brandstubbe'

# Row 5
$ws.Range("A5").Value = 111609174
$ws.Range("B5").Value = 90854
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 2079
$ws.Range("F5").Value = 'Nordtagging'
$ws.Range("G5").Value = 'Odonticium romellii'
$ws.Range("H5").Value = '(S.Lundell) Parmasto'
$ws.Range("I5").Value = ""
$ws.Range("Q5").Value = 514788.8674634451
$ws.Range("R5").Value = 6925250.666874606
$ws.Range("AO5").Value = 'mossig silverlåga av tall'

# Row 6
$ws.Range("A6").Value = 111609167
$ws.Range("B6").Value = 77186
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 353
$ws.Range("F6").Value = 'Dvärgbägarlav'
$ws.Range("G6").Value = 'Cladonia parasitica'
$ws.Range("H6").Value = '(Hoffm.) Hoffm.'
$ws.Range("I6").Value = ""
$ws.Range("Q6").Value = 515051.1877758073
$ws.Range("R6").Value = 6925144.938876954
$ws.Range("AO6").Value = 'silverlåga av tall'

# Row 7
$ws.Range("A7").Value = 111609170
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = 'Knärot'
$ws.Range("G7").Value = 'Goodyera repens'
$ws.Range("H7").Value = '(L.) R. Br.'
$ws.Range("I7").Value = '''3'
$ws.Range("Q7").Value = 515035.9338400747
$ws.Range("R7").Value = 6925238.814452391
$ws.Range("AO7").Value = ""

# Row 8
$ws.Range("A8").Value = 111609169
$ws.Range("B8").Value = 96348
$ws.Range("D8").Value = 'VU'
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = 'Knärot'
$ws.Range("G8").Value = 'Goodyera repens'
$ws.Range("H8").Value = '(L.) R. Br.'
$ws.Range("I8").Value = '''4'
$ws.Range("Q8").Value = 515078.8479096842
$ws.Range("R8").Value = 6925177.45879681
$ws.Range("AO8").Value = ""

# Row 9
$ws.Range("A9").Value = 111609176
$ws.Range("B9").Value = 96348
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = 'Knärot'
$ws.Range("G9").Value = 'Goodyera repens'
$ws.Range("H9").Value = '(L.) R. Br.'
$ws.Range("I9").Value = '''30'
$ws.Range("Q9").Value = 514875.4249116365
$ws.Range("R9").Value = 6924962.732657854
$ws.Range("AO9").Value = ""

# Row 10
$ws.Range("A10").Value = 111609168
$ws.Range("B10").Value = 77597
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 864
$ws.Range("F10").Value = 'Knottrig blåslav'
$ws.Range("G10").Value = 'Hypogymnia bitteri'
$ws.Range("H10").Value = '(Lynge) Ahti'
$ws.Range("I10").Value = ""
$ws.Range("Q10").Value = 515085.0087401169
$ws.Range("R10").Value = 6925147.4056778
$ws.Range("AO10").Value = 'tall'
